# Appointment_Outcomes.xlsx - pharmacist dispensed prescriptions for
# appointment 2 (Patient P1001 / Doctor D002) and appointment 3
# (Patient P1002 / Doctor D001): their "Prescription Status" moves
# from PENDING to DISPENSED.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column G is "Prescription Status"; rows 3 and 4 hold appointments 2 and 3.
$ws.Range("G3").Value = "DISPENSED"
$ws.Range("G4").Value = "DISPENSED"
